# Weekly update: a new week of price data is inserted at the top of the
# "Vega Monumental Concepción - Cebolla" block (rows 263-264), pushing all
# subsequent rows down by two. The sheet's used range therefore grows from
# A1:R359 to A1:R361.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new, blank rows at 263-264; everything below shifts down by 2.
$ws.Rows("263:264").Insert()

# --- New row 263 ---
$ws.Range("A263").Value2 = 11
$ws.Range("B263").Value2 = "Vega Monumental Concepción"
$ws.Range("C263").Value2 = "Bíobío"
$ws.Range("D263").Value2 = 44588
$ws.Range("E263").Value2 = 8
$ws.Range("F263").Value2 = 100112004
$ws.Range("G263").Value2 = "Cebolla"
$ws.Range("H263").Value2 = "Sin especificar"
$ws.Range("I263").Value2 = "1a (cosecha)"
$ws.Range("J263").Value2 = 350
$ws.Range("K263").Value2 = 4500
$ws.Range("L263").Value2 = 5000
$ws.Range("M263").Value2 = 4714
$ws.Range("N263").Value2 = "$/malla 18 kilos"
$ws.Range("O263").Value2 = "Región de O'Higgins"
$ws.Range("P263").Value2 = 262
$ws.Range("Q263").Value2 = 18
$ws.Range("R263").Value2 = "Hortaliza"

# --- New row 264 ---
$ws.Range("A264").Value2 = 11
$ws.Range("B264").Value2 = "Vega Monumental Concepción"
$ws.Range("C264").Value2 = "Bíobío"
$ws.Range("D264").Value2 = 44588
$ws.Range("E264").Value2 = 8
$ws.Range("F264").Value2 = 100112004
$ws.Range("G264").Value2 = "Cebolla"
$ws.Range("H264").Value2 = "Sin especificar"
$ws.Range("I264").Value2 = "2a (cosecha)"
$ws.Range("J264").Value2 = 300
$ws.Range("K264").Value2 = 4000
$ws.Range("L264").Value2 = 4000
$ws.Range("M264").Value2 = 4000
$ws.Range("N264").Value2 = "$/malla 18 kilos"
$ws.Range("O264").Value2 = "Región de O'Higgins"
$ws.Range("P264").Value2 = 222
$ws.Range("Q264").Value2 = 18
$ws.Range("R264").Value2 = "Hortaliza"

# Make sure the date cells keep the same date/time number format used by
# the rest of column D (style index 2 in the original workbook).
$ws.Range("D263:D264").NumberFormat = $ws.Range("D265").NumberFormat
